# logboek.xlsx - "gitignore update, prototype 1, logboek"
#
# For each of the three logbook sheets (Gerben, Wester, Garon) a new day
# entry (44886 = the next work day) is appended after the existing second
# entry. That pushes what used to be the "last row" formatting (open
# bottom-border box) down to the new row, and turns the former
# second-day box (rows 6-7) into a normal, closed 2-row box - matching
# the box drawn for the very first day (rows 3-5).
#
# We reproduce the exact look by first snapshotting the live formatting
# of rows 3 (top of box), 5 (bottom of box) and 6 (the "open" box style)
# with Copy / PasteSpecial (formats only) BEFORE we overwrite anything,
# then writing in the new values.

$wb = $excel.ActiveWorkbook

function Set-DayEntry($ws, $activity) {

    # --- snapshot current formatting of the template rows we need ---
    # row 6 currently carries the "open ended" box formatting that the
    # new row 8 must inherit.
    $ws.Range("B6:F6").Copy()
    $ws.Range("B8:F8").PasteSpecial(-4122)   # xlPasteFormats

    # row 3 is the top of a closed box -> becomes the new look of row 6
    $ws.Range("B3:F3").Copy()
    $ws.Range("B6:F6").PasteSpecial(-4122)   # xlPasteFormats

    # row 5 is the bottom of a closed box -> becomes the new look of row 7
    # (column C keeps its own distinct "activity" box style, untouched)
    $ws.Range("B5").Copy()
    $ws.Range("B7").PasteSpecial(-4122)      # xlPasteFormats
    $ws.Range("D5:F5").Copy()
    $ws.Range("D7:F7").PasteSpecial(-4122)   # xlPasteFormats

    $ws.Application.CutCopyMode = $false

    # --- values for the freshly appended day row ---
    $ws.Range("B8").Value = 44886
    $ws.Range("C8").Value = $activity
    $ws.Range("D8").Value = 0.375
    $ws.Range("E8").Value = 0.4826388888888889
    $ws.Range("F8").Value = 0.1111111111111111
}

$wsGerben = $wb.Worksheets.Item("Gerben")
$wsWester = $wb.Worksheets.Item("Wester")
$wsGaron  = $wb.Worksheets.Item("Garon")

Set-DayEntry $wsGerben "Werken aan home page prototype 1"
Set-DayEntry $wsWester "Gewerkt aan prototype 3"
Set-DayEntry $wsGaron  "Gewerkt aan prototype 2"

# Wester and Garon's worked-until time differs slightly from Gerben's.
$wsWester.Range("E8").Value = 0.4861111111111111
$wsGaron.Range("E8").Value = 0.4861111111111111

# --- selections / active sheet, last touched = Garon ---
$wsGerben.Activate()
$wsGerben.Range("C13").Select()

$wsWester.Activate()
$wsWester.Range("I10").Select()

$wsGaron.Activate()
$wsGaron.Range("C12").Select()

Write-Output "logboek.xlsx updated"
